$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 34; existing rows 34-84 shift down to 35-85.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new data record.
$ws.Cells.Item(34, 1).Value = 10
$ws.Cells.Item(34, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(34, 3).Value = "La Araucanía"
$ws.Cells.Item(34, 4).Value = 44477
$ws.Cells.Item(34, 5).Value = 9
$ws.Cells.Item(34, 6).Value = 100112012
$ws.Cells.Item(34, 7).Value = "Espinaca"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 20
$ws.Cells.Item(34, 11).Value = 11000
$ws.Cells.Item(34, 12).Value = 11000
$ws.Cells.Item(34, 13).Value = 11000
$ws.Cells.Item(34, 14).Value = "$/docena de atados"
$ws.Cells.Item(34, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(34, 16).Value = 3667
$ws.Cells.Item(34, 17).Value = 3
$ws.Cells.Item(34, 18).Value = "Hortaliza"
